$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$NL = [char]10

# --- Row 15 / Row 16: new content (B15, C15, C16, D16, D15 in this exact
#     order so the shared-string table is populated the same way the
#     original author's session produced it) ---
$ws.Range("B15").Value = "Interaction Mode Combo Box" + $NL + "Manual Control Focus button" + $NL + "Configuration Controls" + $NL + "Greenhouse data controls"
$ws.Range("B15").WrapText = $true
$ws.Range("B15").VerticalAlignment = -4160

$ws.Range("C15").Value = "Tests whether appropriate controls will enable/disable when an interaction mode is selected."
$ws.Range("C15").WrapText = $true
$ws.Range("C15").VerticalAlignment = -4160

# --- Row 16: brand new row (previously an empty gap between row 15 and row 17) ---
$ws.Range("A16").Value = $ws.Range("A15").Text
$ws.Range("A16").WrapText = $true
$ws.Range("A16").VerticalAlignment = -4160

$ws.Range("B16").Value = "N/A"
$ws.Range("B16").WrapText = $true
$ws.Range("B16").VerticalAlignment = -4160

$ws.Range("C16").Value = "Tests whether the buggy goes into autonomous mode when the mode is selected and stops when another mode is selected"
$ws.Range("C16").WrapText = $true
$ws.Range("C16").VerticalAlignment = -4160

$ws.Range("D16").Value = "Buggy goes into autonomous mode when the mode is selected and stops working in autonomous mode when another interaction mode is selected."
$ws.Range("D16").WrapText = $true

$ws.Rows.Item(16).RowHeight = 60

$ws.Range("D15").Value = "Manual Mode:" + $NL + "- Control Focus Button enabled" + $NL + "- Controls under greenhouse data section on the gui enabled." + $NL + "Configuration Mode:" + $NL + "- Configuration controls enabled" + $NL + "Autonomous Mode:" + $NL + "None of the Manual, Configuration or Greenhouse data controls are enabled" + $NL
$ws.Range("D15").WrapText = $true

$ws.Rows.Item(15).RowHeight = 150

# --- Row 17: gains a new A value; existing E17 cell is untouched ---
$ws.Range("A17").Value = "Controlling buggy in Manual Mode"
$ws.Range("A17").WrapText = $true
$ws.Range("A17").VerticalAlignment = -4160

$ws.Rows.Item(17).RowHeight = 30

# --- View state: scrolled so row 13 is at top, selection on B17 ---
$excel.ActiveWindow.ScrollRow = 13
$ws.Range("B17").Select()
